$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 20190717 (row 18) - fill in the remaining checkboxes; cleanliness (M) is an
# exception marked with a cross instead of a check.
$ws.Range("D18").Value = "√"
$ws.Range("G18").Value = "√"
$ws.Range("H18").Value = "√"
$ws.Range("L18").Value = "√"
$ws.Range("M18").Value = "×"
$ws.Range("N18").Value = "√"

# 20190718 (row 19) - start a new day, only the first two items done so far.
$ws.Range("A19").Value = 20190718
$ws.Range("B19").Value = "√"
$ws.Range("O19").Value = "√"

$ws.Range("M19").Select() | Out-Null
